$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at row 70 (pushes existing rows 70..101 down to 71..102)
$ws.Rows.Item(70).Insert()

# Populate the new row 70 with the new weekly price record
$ws.Cells.Item(70, 1).Value = 8
$ws.Cells.Item(70, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(70, 3).Value = "Coquimbo"
$ws.Cells.Item(70, 4).Value = 44917
$ws.Cells.Item(70, 5).Value = 4
$ws.Cells.Item(70, 6).Value = 100112030
$ws.Cells.Item(70, 7).Value = "Poroto granado"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 520
$ws.Cells.Item(70, 11).Value = 34000
$ws.Cells.Item(70, 12).Value = 35000
$ws.Cells.Item(70, 13).Value = 34500
$ws.Cells.Item(70, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(70, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(70, 16).Value = 1380
$ws.Cells.Item(70, 17).Value = 25
$ws.Cells.Item(70, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same number format as the other date cells in column D
$ws.Cells.Item(70, 4).NumberFormat = $ws.Cells.Item(71, 4).NumberFormat
